# "Merge k Sorted Lists Complete"
#
# The workbook tracks LeetCode Blind 75 problems on the "Blind 75" sheet.
# This change marks "Merge k Sorted Lists" (row 41) as completed, and fixes
# up difficulty ratings that had apparently been entered on the wrong rows:
#   - Course Schedule              (row 32): Hard   -> Medium
#   - Merge Two Sorted Lists       (row 40): Hard   -> Medium
#   - Merge k Sorted Lists         (row 41): Medium -> Hard, now completed
#   - Binary Tree Level Order Trav (row 58): Easy   -> Medium

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blind 75")

# Difficulty corrections
$ws.Range("D32").Value = "Medium"
$ws.Range("D40").Value = "Medium"
$ws.Range("D58").Value = "Medium"

# Merge k Sorted Lists: mark as Hard / completed with notes + runtime
$ws.Range("D41").Value = "Hard"
$ws.Range("E41").Value = "X"
$ws.Range("F41").Value = "Divide and conquer using merge 2 sorted lists"
$ws.Range("G41").Value = "O(Nlog(k))"

# Update the active selection to reflect where the author ended up working
$ws.Range("F33").Select()
